$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.735.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.032.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.49%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.47"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.07%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.442"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.53"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.56%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.368"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.555.70"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.80"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000166"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.747.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.16%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.033.35"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.60%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.10"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.10%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.500"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.35%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.76"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.97%  "

$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.73%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.45%  "

$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0937"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.53%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.79"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.00%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.52"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +8.13%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.81"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.75%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.45%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.77"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.74"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.57%  "

$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.16"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.51%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.89"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.62%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.28"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.68%  "

$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.01"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.14%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0687"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.61%  "

$ws.Range("B39").Value = "RenzoRestakedETH"
$ws.Range("C39").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.068.37"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.51"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.22%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.88"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.23%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.316.09"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.54%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.657"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.14%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.43"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.51%  "

$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.05"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.81%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0241"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.19%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.70"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.86%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.86"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.05%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0895"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.00%  "
